{"js": "// Update the date line and each three-digit-division answer cell.\n// Each (before -> after) pair is unique in the document, so a targeted\n// search + full-text replace keeps existing run formatting intact.\nconst replacements = [\n  [\"2025-02-15 Saturday\", \"2025-02-16 Sunday\"],\n  [\"454\u00f77=64, 6\", \"866\u00f75=173, 1\"],\n  [\"108\u00f73=36, 0\", \"663\u00f73=221, 0\"],\n  [\"858\u00f72=429, 0\", \"868\u00f75=173, 3\"],\n  [\"394\u00f79=43, 7\", \"741\u00f72=370, 1\"],\n  [\"975\u00f72=487, 1\", \"993\u00f75=198, 3\"],\n  [\"598\u00f73=199, 1\", \"714\u00f76=119, 0\"],\n  [\"544\u00f79=60, 4\", \"205\u00f77=29, 2\"],\n  [\"395\u00f79=43, 8\", \"301\u00f78=37, 5\"],\n  [\"919\u00f76=153, 1\", \"547\u00f76=91, 1\"],\n  [\"763\u00f73=254, 1\", \"763\u00f74=190, 3\"],\n  [\"188\u00f73=62, 2\", \"656\u00f76=109, 2\"],\n  [\"219\u00f73=73, 0\", \"745\u00f77=106, 3\"],\n  [\"990\u00f79=110, 0\", \"126\u00f77=18, 0\"],\n  [\"523\u00f73=174, 1\", \"230\u00f73=76, 2\"],\n  [\"631\u00f77=90, 1\", \"743\u00f75=148, 3\"],\n  [\"981\u00f73=327, 0\", \"168\u00f72=84, 0\"],\n  [\"455\u00f79=50, 5\", \"177\u00f74=44, 1\"],\n  [\"311\u00f76=51, 5\", \"388\u00f77=55, 3\"],\n  [\"555\u00f79=61, 6\", \"831\u00f79=92, 3\"],\n  [\"114\u00f75=22, 4\", \"523\u00f75=104, 3\"],\n  [\"775\u00f77=110, 5\", \"291\u00f76=48, 3\"],\n  [\"475\u00f73=158, 1\", \"341\u00f77=48, 5\"],\n  [\"857\u00f76=142, 5\", \"406\u00f79=45, 1\"],\n  [\"322\u00f78=40, 2\", \"394\u00f75=78, 4\"],\n  [\"980\u00f76=163, 2\", \"552\u00f79=61, 3\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date line and each three-digit-division answer cell.\n# Each (old -> new) pair is unique in the document, so Find/Replace\n# targeted at exact text keeps existing run formatting intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-02-15 Saturday\"; New = \"2025-02-16 Sunday\" }\n    @{ Old = \"454\u00f77=64, 6\"; New = \"866\u00f75=173, 1\" }\n    @{ Old = \"108\u00f73=36, 0\"; New = \"663\u00f73=221, 0\" }\n    @{ Old = \"858\u00f72=429, 0\"; New = \"868\u00f75=173, 3\" }\n    @{ Old = \"394\u00f79=43, 7\"; New = \"741\u00f72=370, 1\" }\n    @{ Old = \"975\u00f72=487, 1\"; New = \"993\u00f75=198, 3\" }\n    @{ Old = \"598\u00f73=199, 1\"; New = \"714\u00f76=119, 0\" }\n    @{ Old = \"544\u00f79=60, 4\"; New = \"205\u00f77=29, 2\" }\n    @{ Old = \"395\u00f79=43, 8\"; New = \"301\u00f78=37, 5\" }\n    @{ Old = \"919\u00f76=153, 1\"; New = \"547\u00f76=91, 1\" }\n    @{ Old = \"763\u00f73=254, 1\"; New = \"763\u00f74=190, 3\" }\n    @{ Old = \"188\u00f73=62, 2\"; New = \"656\u00f76=109, 2\" }\n    @{ Old = \"219\u00f73=73, 0\"; New = \"745\u00f77=106, 3\" }\n    @{ Old = \"990\u00f79=110, 0\"; New = \"126\u00f77=18, 0\" }\n    @{ Old = \"523\u00f73=174, 1\"; New = \"230\u00f73=76, 2\" }\n    @{ Old = \"631\u00f77=90, 1\"; New = \"743\u00f75=148, 3\" }\n    @{ Old = \"981\u00f73=327, 0\"; New = \"168\u00f72=84, 0\" }\n    @{ Old = \"455\u00f79=50, 5\"; New = \"177\u00f74=44, 1\" }\n    @{ Old = \"311\u00f76=51, 5\"; New = \"388\u00f77=55, 3\" }\n    @{ Old = \"555\u00f79=61, 6\"; New = \"831\u00f79=92, 3\" }\n    @{ Old = \"114\u00f75=22, 4\"; New = \"523\u00f75=104, 3\" }\n    @{ Old = \"775\u00f77=110, 5\"; New = \"291\u00f76=48, 3\" }\n    @{ Old = \"475\u00f73=158, 1\"; New = \"341\u00f77=48, 5\" }\n    @{ Old = \"857\u00f76=142, 5\"; New = \"406\u00f79=45, 1\" }\n    @{ Old = \"322\u00f78=40, 2\"; New = \"394\u00f75=78, 4\" }\n    @{ Old = \"980\u00f76=163, 2\"; New = \"552\u00f79=61, 3\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n    if (-not $result) {\n        throw \"Text not found: $($r.Old)\"\n    }\n}\n"}
